$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# ---------------------------------------------------------------------------
# 1) Row 211 (xylem_total_pressure): correct the over-inflow values in I:T
#    from -30000 to -100000 so they match the G/H columns.
# ---------------------------------------------------------------------------
$cols211 = @("I","J","K","L","M","N","O","P","Q","R","S","T")
foreach ($col in $cols211) {
    $ws.Range($col + "211").Value2 = -100000
}

# ---------------------------------------------------------------------------
# 2) Insert a new parameter row at row 221 ("diffusion_xylem"), pushing the
#    existing rows 221-224 down to 222-225.
# ---------------------------------------------------------------------------
$ws.Rows.Item(221).Insert()

# Copy formatting for the new row from the row directly below it (the row
# that used to be 221, now shifted to 222) so styles match exactly.
$ws.Range("A222:T222").Copy()
$ws.Range("A221:T221").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row's content.
$ws.Range("A221").Value2 = "diffusion_xylem"
$ws.Range("B221").Value2 = "parameter"
$ws.Range("C221").Value2 = "root_cynaps"
$ws.Range("D221").Value2 = "None"
$ws.Range("E221").Value2 = "Diffusion paramenter for exchanges between xylem apoplasm and cortex symplasm"
$ws.Range("F221").Value2 = "g.s-1.m-3"

# Leave G221:T221 empty (no values), matching the source edit.

# ---------------------------------------------------------------------------
# 3) Update the view: selection moves to I229.
# ---------------------------------------------------------------------------
$ws.Range("I229").Select()
